# Applies the manuscript edits described in the commit "sdy and gradient results".
$d = $word.ActiveDocument

$wmNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"'

function Get-FoundRange([string]$searchText) {
    $probe = $d.Content
    $ok = $probe.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Text not found: $searchText"
    }
    return $d.Range($probe.Start, $probe.End)
}

# Builds the oMath XML for "<number><superscript>th</superscript>" e.g. 2.5th, 97.5th
function OMath-Ordinal([string]$number) {
    return "<m:oMath><m:sSup><m:e><m:r><m:t>$number</m:t></m:r></m:e><m:sup><m:r><m:t>t</m:t></m:r><m:r><m:t>h</m:t></m:r></m:sup></m:sSup></m:oMath>"
}

# --- Edit 1: fix typo "possibilty" -> "possibility" ---------------------------
$d.Content.Find.Execute("possibilty", $true, $false, $false, $false, $false, $true, 1, $false, "possibility", 2) | Out-Null

# --- Edit 2: append "(2.5th to 97.5th percentile)" after "...95 % credible intervals" ---
$old = "in eq. ? and ?. To verify that the modern gradient can be approximated with limited sampling, we resampled modern sea surface temperatures at modern latitudes corresponding the palaeolatitudes of the Eocene samples. This process was repeated 100 times, randomly chosing a longitude for each latitude in each repetition. The parameters of the non-hierarchical temperature model were estimated for each of the 100 samples, and the iterations after burn-in from the posterior of the parameters were pooled to generate the resulting median temperature gradient with 95 % credible intervals."
$prefix = "in eq. ? and ?. To verify that the modern gradient can be approximated with limited sampling, we resampled modern sea surface temperatures at modern latitudes corresponding the palaeolatitudes of the Eocene samples. This process was repeated 100 times, randomly chosing a longitude for each latitude in each repetition. The parameters of the non-hierarchical temperature model were estimated for each of the 100 samples, and the iterations after burn-in from the posterior of the parameters were pooled to generate the resulting median temperature gradient with 95 % credible intervals ("
$r = Get-FoundRange($old)
$xml = '<w:p ' + $wmNs + '>' `
    + '<w:r><w:t xml:space="preserve">' + $prefix + '</w:t></w:r>' `
    + (OMath-Ordinal "2.5") `
    + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve">to</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
    + (OMath-Ordinal "97.5") `
    + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve">percentile).</w:t></w:r>' `
    + '</w:p>'
$r.InsertXML($xml)

# --- Edit 3: "...calculating the median (2.5th, 97.5th percentile) from all pairs of iterations." ---
$old = "Differences between Eocene and modern temperatures at a certain latitude were calculated by randomly pairing all iterations of the posterior from the Eocene and modern temperature gradient model, calculating the Eocene and modern temperature using the respective iterations, taking the difference, and then calculating the median (2.5th, 97.5th percentile) from all pairs of iterations."
$prefix = "Differences between Eocene and modern temperatures at a certain latitude were calculated by randomly pairing all iterations of the posterior from the Eocene and modern temperature gradient model, calculating the Eocene and modern temperature using the respective iterations, taking the difference, and then calculating the median ("
$r = Get-FoundRange($old)
$xml = '<w:p ' + $wmNs + '>' `
    + '<w:r><w:t xml:space="preserve">' + $prefix + '</w:t></w:r>' `
    + (OMath-Ordinal "2.5") `
    + '<w:r><w:t xml:space="preserve">,</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
    + (OMath-Ordinal "97.5") `
    + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve">percentile) from all pairs of iterations.</w:t></w:r>' `
    + '</w:p>'
$r.InsertXML($xml)

# --- Edit 4: "...weighted mean of the median (2.5th, 97.5th percentile) temperature estimates..." ---
$nbsp = [char]0x00A0
$old = "Global average temperatures were calculated by taking the weighted mean of the median (2.5th, 97.5th percentile) temperature estimates in 1 degree latitudinal bins. The weights were set to the proportion of global surface area in each latitudinal bin, i.e." + $nbsp + "decreasing with increasing latitude as:"
$prefix = "Global average temperatures were calculated by taking the weighted mean of the median ("
$suffix = " temperature estimates in 1 degree latitudinal bins. The weights were set to the proportion of global surface area in each latitudinal bin, i.e." + $nbsp + "decreasing with increasing latitude as:"
$r = Get-FoundRange($old)
$xml = '<w:p ' + $wmNs + '>' `
    + '<w:r><w:t xml:space="preserve">' + $prefix + '</w:t></w:r>' `
    + (OMath-Ordinal "2.5") `
    + '<w:r><w:t xml:space="preserve">,</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
    + (OMath-Ordinal "97.5") `
    + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve">percentile)' + $suffix + '</w:t></w:r>' `
    + '</w:p>'
$r.InsertXML($xml)

# --- Edit 5: split the Eocene-gradient paragraph into two, fix "modelled"->"modeled", add sdy sentence ---
$old = "The modelled Eocene temperature gradient is starkly different from the modern (Fig 4). Modelled, median equatorial temperatures are 4.2 (95% CI: 0.2 - 8.3) degC higher for the Eocene, and polar temperatures are higher by 25.0 (95% CI: 17.0 - 29.1) degC. This results in a strongly flattened latitudinal temperature gradient of 9.0 (95% CI: 2.5 - 17.8) degC for the Eocene, as opposed to 29.6 degC for the modern. The high variability of early Eocene palaeotemperature proxies, particularly in the mid-latitudes, and the scarcity of high-latitude data, results in substantial uncertainties in the modelled temperature gradient, but " + [char]0x2026
$r = Get-FoundRange($old)
$para = $r.Paragraphs(1)

# Replace the (whole) paragraph's content, keeping its "FirstParagraph" style, with the
# shortened first sentence.
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
    + '<w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr>' `
    + '<w:r><w:t xml:space="preserve">The modelled Eocene temperature gradient is starkly different from the modern (Fig 4). Modelled, median equatorial temperatures are 4.2 (95% CI: 0.2 - 8.3) degC higher for the Eocene, and polar temperatures are higher by 25.0 (95% CI: 17.0 - 29.1) degC. This results in a strongly flattened latitudinal temperature gradient of 9.0 (95% CI: 2.5 - 17.8) degC for the Eocene, as opposed to 29.6 degC for the modern.</w:t></w:r>' `
    + '</w:p>'
$para.Range.InsertXML($xml1)

# Insert a new "BodyText"-styled paragraph straight after it with the new sdy sentence.
$endOfPara = $para.Range
$endOfPara.Collapse(0)
$endOfPara.InsertParagraphAfter()
$newPara = $para.Next()
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
    + '<w:pPr><w:pStyle w:val="BodyText"/></w:pPr>' `
    + '<w:r><w:t xml:space="preserve">The high variability of early Eocene palaeotemperature proxies, particularly in the mid-latitudes, and the scarcity of high-latitude data, results in substantial uncertainties in the modeled temperature gradient. This is reflected in the residual standard deviation, which is much higher for the early Eocene gradient, 4.9 (95% CI: 3.8 - 6.5) degC, than for the modern gradient, resampled at early Eocene latitudes, 2.2 (95% CI: 1.6 - 3.1) degC.</w:t></w:r>' `
    + '</w:p>'
$newPara.Range.InsertXML($xml2)

Write-Output "done"
